$d = $word.ActiveDocument

# Update the date heading (first paragraph)
$d.Paragraphs.Item(1).Range.Text = "2024-01-11 Thursday"

# Update the multiplication table. The table has 20 rows x 5 columns,
# with actual content only in rows 1, 5, 10, 15, 20 (the rest are
# blank spacer rows). Addressing cells positionally avoids any
# ambiguity from duplicate text values (e.g. "42x77=3234" appears
# twice in the original, and "48x96=4608" appears as both a source
# and a target value).
$t = $d.Tables.Item(1)

$newValues = @{
    1  = @("32×15=480", "48×96=4608", "56×78=4368", "60×35=2100", "87×81=7047")
    5  = @("91×98=8918", "86×89=7654", "99×88=8712", "83×58=4814", "49×80=3920")
    10 = @("25×57=1425", "50×66=3300", "95×72=6840", "76×92=6992", "90×98=8820")
    15 = @("31×84=2604", "38×63=2394", "49×96=4704", "48×57=2736", "24×27=648")
    20 = @("34×38=1292", "33×21=693", "92×95=8740", "23×42=966", "64×15=960")
}

foreach ($rowIndex in $newValues.Keys) {
    $cols = $newValues[$rowIndex]
    for ($c = 1; $c -le $cols.Length; $c++) {
        $t.Cell($rowIndex, $c).Range.Text = $cols[$c - 1]
    }
}

Write-Host "Edit complete"
